$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (rows 3-9) ---
# C5: PO Date        15-Sep-2021 (44454) -> 13-Dec-2021 (44543)
$ws.Range("C5").Value = 44543
# C8: Delivery Date  20-Sep-2021 (44459) -> 13-Dec-2021 (44543)
$ws.Range("C8").Value = 44543

# --- Detail line 1 (row 13) ---
# E13: Quantity 1200 -> 10000
$ws.Range("E13").Value = 10000
# H13: Note "Note for line 1" -> "Luu y hang de vo"
$ws.Range("H13").Value = "Luu y hang de vo"

# --- Remaining per-line Note cells (H14:H30) are blanked out ---
$ws.Range("H14:H30").ClearContents()

# --- Sheet view: scroll position, zoom level and active selection ---
$window = $excel.ActiveWindow
$window.ScrollRow = 7
$window.ScrollColumn = 1
$window.Zoom = 145
$ws.Range("G17").Select()
